# Sort episodes by Published Date
# Renumber the index column (A) sequentially for the rows that were
# previously offset (rows 4-16 held values 1392-1404 from a second,
# differently-indexed source list). After sorting by Published Date,
# the whole list is treated as one contiguous sequence, so those rows
# are renumbered to continue directly after the existing 0,1 index
# values already present in rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 4; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
